# Second round TouchGFX text updates — Translation sheet (Table8, B3:I799).
#
# Row 3 is the table's header row; it is rewritten to hold the upper-case
# column captions, loses its explicit cell style (s="1"), and drops the
# last three language columns (G3:I3 are removed entirely).
#
# Row 4 becomes the first data row of the table: a single-use text id
# ("SingleUseId1"), the "Default" typography, "Center" alignment, the
# "Hello" translation for the GB (English) column, and "LTR" direction.
#
# NOTE: writing directly into a ListObject's header row normally makes
# Excel resync the table's column names to match the new cell text. The
# authored workbook keeps the table's column definitions (Table8 in
# xl/tables/table3.xml) untouched, so the table is temporarily resized
# away from row 3 while the header cells are edited, then resized back.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")
$lo = $ws.ListObjects.Item("Table8")
$originalRange = $lo.Range.Address()

# Move the table off row 3 so editing the header cells doesn't resync
# the ListObject's column names.
$lo.Resize($ws.Range("B810:I811"))

# --- Row 3: header row -----------------------------------------------
$ws.Range("B3").Value = "TEXT ID"
$ws.Range("C3").Value = "TYPOGRAPHY NAME"
$ws.Range("D3").Value = "ALIGNMENT"
$ws.Range("E3").Value = "GB"
$ws.Range("F3").Value = "DIRECTION"

# Drop the old per-cell style (s="1") from the header row.
$ws.Range("B3:F3").Style = "Normal"

# The trailing language columns are no longer populated on this row.
$ws.Range("G3:I3").Clear()

# Restore the table back to its original extent.
$lo.Resize($ws.Range($originalRange.Replace("$", "")))

# --- Row 4: first data row --------------------------------------------
$ws.Range("B4").Value = "SingleUseId1"
$ws.Range("C4").Value = "Default"
$ws.Range("D4").Value = "Center"
$ws.Range("E4").Value = "Hello"
$ws.Range("F4").Value = "LTR"

# New cells inherit the column default style (s="1"); strip it so the
# row matches the unstyled data rows produced by the tool.
$ws.Range("B4:F4").Style = "Normal"
